$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "R0" column (C) register trace from 64 to 40 for rows 4-35 ---
$ws.Range("C4:C35").Value = 40

# --- Update text labels that mention the expected value (64 -> 45) ---
$ws.Range("B39").Value = "45, 0"
$ws.Range("M35").Value = "Test succeedes if 45 written in addr 65"
$ws.Range("L39").Value = "45 if test passes, 0 if test fails"

# --- Column F (R3): replace hard-coded 59 with formula C+E, shared over F9:F35 ---
$ws.Range("F8").Formula = "=C8+E8"
$ws.Range("F9:F35").Formula = "=C9+E9"

# --- Column G (R4): update 64 -> 45. Rows 9 & 10 get explicit formulas, rest get literal values ---
$ws.Range("G9").Formula = "=C9-E9"
$ws.Range("G10").Formula = "=C10-E10"

$ws.Range("G12").Value = 45
$ws.Range("G13").Value = 45
$ws.Range("G15").Value = 45
$ws.Range("G17").Value = 45
$ws.Range("G18").Value = 45
$ws.Range("G19").Value = 45
$ws.Range("G21").Value = 45
$ws.Range("G23").Value = 45
$ws.Range("G24").Value = 45
$ws.Range("G26").Value = 45
$ws.Range("G27").Value = 45
$ws.Range("G28").Value = 45
$ws.Range("G30").Value = 45
$ws.Range("G31").Value = 45
$ws.Range("G32").Value = 45
$ws.Range("G33").Value = 45
$ws.Range("G34").Value = 45
$ws.Range("G35").Value = 45

# --- Recalculate so cached formula results match ---
$excel.Calculate()

# --- Update selection to match the final cursor position recorded in the file ---
$ws.Range("L40").Select()
